$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 71 (year 2022) values
$ws.Range("B71").Value = 98
$ws.Range("C71").Value = 10

# Add new row 72 (year 2023)
$ws.Range("A72").Value = 2023
$ws.Range("B72").Value = 10
$ws.Range("C72").Value = 6
